$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, shifting existing rows 198:217 down to 199:218.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly price record.
$ws.Range("A198").Value = 4
$ws.Range("B198").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C198").Value = "Los Lagos"
$ws.Range("D198").Value = 44578
$ws.Range("E198").Value = 10
$ws.Range("F198").Value = 100112040
$ws.Range("G198").Value = "Cilantro"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 40
$ws.Range("K198").Value = 10000
$ws.Range("L198").Value = 10000
$ws.Range("M198").Value = 10000
$ws.Range("N198").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O198").Value = "Región de La Araucanía"
$ws.Range("P198").Value = 5000
$ws.Range("Q198").Value = 2
$ws.Range("R198").Value = "Hortaliza"
